$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.5
$ws.Range("AK2").Value = 41
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 10
